$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.913.07"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "2.294.36"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D5").Value = "'300.70"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "'98.89"
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("D7").Value = "'0.504"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +3.14%  "
$ws.Range("D10").Value = "'36.09"
$ws.Range("E10").Value = "  +7.60%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "'18.40"
$ws.Range("E12").Value = "  +9.48%  "
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "'6.95"
$ws.Range("E14").Value = "  +2.70%  "
$ws.Range("D15").Value = "2.650.42"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "2.572.37"
$ws.Range("E16").Value = "  +11.87%  "
$ws.Range("D17").Value = "'0.800"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "42.824.25"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "'12.57"
$ws.Range("E19").Value = "  +8.44%  "
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("D22").Value = "'67.67"
$ws.Range("D23").Value = "'235.49"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "'2.23"
$ws.Range("E24").Value = "  +10.94%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").Value = "'25.01"
$ws.Range("E27").Value = "  +2.65%  "
$ws.Range("D28").Value = "'2.35"
$ws.Range("E28").Value = "  +14.65%  "
$ws.Range("E29").Value = "  +1.74%  "
$ws.Range("D30").Value = "'167.22"
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").Value = "'9.12"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "'17.65"
$ws.Range("E34").Value = "  +5.10%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "'4.65"
$ws.Range("E35").Value = "  -1.51%  "
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("E39").Value = "  +2.25%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").Value = "'0.0291"
$ws.Range("E43").Value = "  +3.82%  "
$ws.Range("D44").Value = "1.976.84"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").Value = "'10.10"
$ws.Range("E45").Value = "  +3.03%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.88"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'17.48"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").Value = "'55.14"
$ws.Range("E48").Value = "  +4.56%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.54"
$ws.Range("E49").Value = "  +3.40%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.517.90"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").Value = "'70.72"
$ws.Range("E51").Value = "  +1.39%  "
